$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.336
$ws.Range("A3").Value = -21.445
$ws.Range("B5").Value = 6.712999999999999
$ws.Range("D5").Value = -8.395
$ws.Range("E7").Value = 13.078
$ws.Range("D9").Value = -7.721000000000001
$ws.Range("D11").Value = -8.171000000000001
$ws.Range("E11").Value = 12.931
$ws.Range("A14").Value = -20.779
$ws.Range("E19").Value = 12.931
$ws.Range("A21").Value = -20.733
$ws.Range("D21").Value = -7.911
$ws.Range("E21").Value = 13.153
$ws.Range("A23").Value = -21.584
$ws.Range("A25").Value = -22.078
